# Generate Report for Archive
# Updates the localization-status report: the hand-off stage moves from
# "Ready for handoff" to "In Translation" on every sheet that tracks it,
# and the Status column is re-sized to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells (columns E & F, row 2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status column (column C, row 2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
